$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two matches that share the same kickoff date/time got recorded with
# their rows swapped (everything except the row's running index in column A).
# Fix each pair by swapping columns B:AD between the two rows.

$pairs = @(
    @(115, 116),
    @(165, 166),
    @(193, 194),
    @(268, 269),
    @(271, 272)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}

# The row-271/272 pair also carries two re-priced PL figures (PLA / PL_Aha)
# that weren't simple carry-overs from the swap, so set them explicitly to
# match the refreshed base.
$ws.Range("Z272").Value2 = 1.45
$ws.Range("AB272").Value2 = -1
